# Auto-generated script to apply scheduled market-data refresh to Bahamut_Profits workbook
# Updates currentAveragePrice / Leve profit columns (H-N) across ALC, ARM, BSM, CRP, GSM, LTW, WVR sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1801.6316
$ws.Range("I98").Value = 1920.3125
$ws.Range("K98").Value = 1920.3125
$ws.Range("M98").Value = -422.3125

$ws.Range("H113").Value = 3208.7693
$ws.Range("J113").Value = 2856.6667
$ws.Range("L113").Value = 2856.6667
$ws.Range("N113").Value = -9364.6667

$ws.Range("H122").Value = 1801.6316
$ws.Range("I122").Value = 1920.3125
$ws.Range("K122").Value = 5760.9375
$ws.Range("M122").Value = -3310.9375

$ws.Range("H129").Value = 772583.25
$ws.Range("I129").Value = 354.45456
$ws.Range("J129").Value = 1002164.75
$ws.Range("K129").Value = 1063.36368
$ws.Range("L129").Value = 3006494.25
$ws.Range("M129").Value = 3936.63632
$ws.Range("N129").Value = -3016494.25

$ws.Range("H132").Value = 1525.279
$ws.Range("I132").Value = 1525.279
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4575.837
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -2045.837
$ws.Range("N132").ClearContents()

$ws.Range("H137").Value = 1423.1111
$ws.Range("I137").Value = 1120.4
$ws.Range("J137").Value = 1801.5
$ws.Range("K137").Value = 3361.2
$ws.Range("L137").Value = 5404.5
$ws.Range("M137").Value = -811.2000000000003
$ws.Range("N137").Value = -10504.5

$ws.Range("H141").Value = 4752.5
$ws.Range("I141").Value = 4502.857
$ws.Range("J141").Value = 6500
$ws.Range("K141").Value = 13508.571
$ws.Range("L141").Value = 19500
$ws.Range("M141").Value = -8328.571
$ws.Range("N141").Value = -29860

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 23755.6
$ws.Range("I2").Value = 29481.629
$ws.Range("J2").Value = 3714.5
$ws.Range("K2").Value = 29481.629
$ws.Range("L2").Value = 3714.5
$ws.Range("M2").Value = -29368.629
$ws.Range("N2").Value = -3940.5

$ws.Range("H31").Value = 8000
$ws.Range("I31").Value = 8000
$ws.Range("K31").Value = 8000
$ws.Range("M31").Value = -7706

$ws.Range("H116").Value = 23755.6
$ws.Range("I116").Value = 29481.629
$ws.Range("J116").Value = 3714.5
$ws.Range("K116").Value = 29481.629
$ws.Range("L116").Value = 3714.5
$ws.Range("M116").Value = -27187.629
$ws.Range("N116").Value = -8302.5

$ws.Range("H132").Value = 1724.5454
$ws.Range("I132").Value = 1419.1538
$ws.Range("K132").Value = 4257.4614
$ws.Range("M132").Value = -1727.4614

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 23755.6
$ws.Range("I3").Value = 29481.629
$ws.Range("J3").Value = 3714.5
$ws.Range("K3").Value = 29481.629
$ws.Range("L3").Value = 3714.5
$ws.Range("M3").Value = -29367.629
$ws.Range("N3").Value = -3942.5

$ws.Range("H23").Value = 49211.2
$ws.Range("I23").Value = 6000
$ws.Range("K23").Value = 6000
$ws.Range("M23").Value = -5717

$ws.Range("H134").Value = 38729.57
$ws.Range("I134").Value = 2493.5908
$ws.Range("K134").Value = 7480.7724
$ws.Range("M134").Value = -4945.7724

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 62867.824
$ws.Range("I31").Value = 4429.4165
$ws.Range("J31").Value = 203120
$ws.Range("K31").Value = 4429.4165
$ws.Range("L31").Value = 203120
$ws.Range("M31").Value = -4134.4165
$ws.Range("N31").Value = -203710

$ws.Range("H34").Value = 62867.824
$ws.Range("I34").Value = 4429.4165
$ws.Range("J34").Value = 203120
$ws.Range("K34").Value = 4429.4165
$ws.Range("L34").Value = 203120
$ws.Range("M34").Value = -4227.4165
$ws.Range("N34").Value = -203524

$ws.Range("H94").Value = 9307.200000000001
$ws.Range("J94").Value = 9451.111000000001
$ws.Range("L94").Value = 9451.111000000001
$ws.Range("N94").Value = -10353.111

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2195.5908
$ws.Range("I126").Value = 2530.3572
$ws.Range("K126").Value = 7591.071599999999
$ws.Range("M126").Value = -5121.071599999999

$ws.Range("H132").Value = 3216.2856
$ws.Range("I132").Value = 2921.8572
$ws.Range("K132").Value = 8765.571599999999
$ws.Range("M132").Value = -6235.571599999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1686.5264
$ws.Range("I7").Value = 1602.125
$ws.Range("J7").Value = 2136.6667
$ws.Range("K7").Value = 1602.125
$ws.Range("L7").Value = 2136.6667
$ws.Range("M7").Value = -1490.125
$ws.Range("N7").Value = -2360.6667

$ws.Range("H40").Value = 723222.0600000001
$ws.Range("I40").Value = 1263926.1
$ws.Range("K40").Value = 1263926.1
$ws.Range("M40").Value = -1263790.1

$ws.Range("H98").Value = 15177.5
$ws.Range("J98").Value = 15177.5
$ws.Range("L98").Value = 15177.5
$ws.Range("N98").Value = -21167.5

$ws.Range("H122").Value = 2382.85
$ws.Range("I122").Value = 2353.25
$ws.Range("J122").Value = 2501.25
$ws.Range("K122").Value = 7059.75
$ws.Range("L122").Value = 7503.75
$ws.Range("M122").Value = -4609.75
$ws.Range("N122").Value = -12403.75

$ws.Range("H126").Value = 1686.5264
$ws.Range("I126").Value = 1602.125
$ws.Range("J126").Value = 2136.6667
$ws.Range("K126").Value = 4806.375
$ws.Range("L126").Value = 6410.000100000001
$ws.Range("M126").Value = -2336.375
$ws.Range("N126").Value = -11350.0001

$ws.Range("H132").Value = 2078.152
$ws.Range("I132").Value = 1972.6364
$ws.Range("J132").Value = 2346
$ws.Range("K132").Value = 5917.9092
$ws.Range("L132").Value = 7038
$ws.Range("M132").Value = -3387.9092
$ws.Range("N132").Value = -12098

$ws.Range("H136").Value = 3232.1304
$ws.Range("I136").Value = 2009.5428
$ws.Range("J136").Value = 7122.1816
$ws.Range("K136").Value = 6028.6284
$ws.Range("L136").Value = 21366.5448
$ws.Range("M136").Value = -3478.6284
$ws.Range("N136").Value = -26466.5448

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 700.8788
$ws.Range("I126").Value = 650.9666999999999
$ws.Range("J126").Value = 1200
$ws.Range("K126").Value = 1952.9001
$ws.Range("L126").Value = 3600
$ws.Range("M126").Value = 517.0999000000002
$ws.Range("N126").Value = -8540

$ws.Range("H132").Value = 1091.9062
$ws.Range("I132").Value = 937.96155
$ws.Range("K132").Value = 2813.88465
$ws.Range("M132").Value = -283.88465

$ws.Range("H136").Value = 1273.9025
$ws.Range("I136").Value = 1163.4857
$ws.Range("J136").Value = 1918
$ws.Range("K136").Value = 3490.4571
$ws.Range("L136").Value = 5754
$ws.Range("M136").Value = -940.4570999999996
$ws.Range("N136").Value = -10854

Write-Output "Applied scheduled market-data refresh to 29 rows across 7 sheets."
